$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "preparer" columns (B = s1cDNAPreparer, E = s2cDNAPreparer) were
# corrected from "BROWN" to "H.BROWN" for every sample row (2-29), to line
# up with the naming convention used by the other metadata template sheets.
for ($r = 2; $r -le 29; $r++) {
    $ws.Cells.Item($r, 2).Value = "H.BROWN"
    $ws.Cells.Item($r, 5).Value = "H.BROWN"
}

# Row 2 had been explicitly styled (Arial 10 / theme text color) while every
# other data row used the sheet's default formatting; restore the default
# "Normal" style there so the whole column is consistent again.
$ws.Range("B2").Style = "Normal"
$ws.Range("E2").Style = "Normal"

# Reflect where the author had scrolled/selected when the file was last
# saved.
$ws.Range("E28:E29").Select()
